$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CastAimedAbility")

# Insert a new column before the existing "Note" column (H) and give it the
# "SearchUnit" header; this shifts the old "Note" header from H1 to I1.
$ws.Range("H1").EntireColumn.Insert()
$ws.Range("H1").Value = "SearchUnit"

# Match the post-edit selection state recorded for this sheet.
$ws.Activate()
$ws.Range("M4").Select()
